$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(1).Insert()
$hl1 = $ws1.Hyperlinks.Item(1)
$hl1.Range = $ws1.Range("B2")
